$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.285.98"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "1.608.59"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.90"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.01"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.486"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0613"
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.15"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").Value = "1.610.70"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.02"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.515"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").Value = "26.321.25"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.90"
$ws.Range("E17").Value = "  +2.15%  "
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.73"
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.27"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.30"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.01"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.85"
$ws.Range("E24").Value = "  +1.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.65"
$ws.Range("E25").Value = "  +1.91%  "
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.121"
$ws.Range("E27").Value = "  -3.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.20"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.56"
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0491"
$ws.Range("E30").Value = "  +4.57%  "
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.19"
$ws.Range("E32").Value = "  +2.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.92"
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.43"
$ws.Range("E34").Value = "  +2.87%  "
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("D36").Value = "1.165.52"
$ws.Range("E36").Value = "  +5.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0166"
$ws.Range("E37").Value = "  +2.32%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.33"
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.787"
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.499"
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.781"
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.25"
$ws.Range("E43").Value = "  +3.29%  "
$ws.Range("D44").Value = "1.751.15"
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.81"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.52"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.09"
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₇0987"
$ws.Range("E49").Value = "  -8.11%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.408"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("E51").Value = "  -0.04%  "
